$wb = $excel.ActiveWorkbook

# Update the "Load-day" sheet: multiply the Load values (B2:B61) by 40
$ws = $wb.Worksheets.Item("Load-day")
for ($i = 2; $i -le 61; $i++) {
    $cell = $ws.Cells.Item($i, 2)
    $cell.Value = $cell.Value() * 40
}

# Make "Load-day" the active/selected sheet (was "Gen-bus")
$ws.Activate()
